$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new data-driven-test rows
$ws.Range("A4").Value = "."
$ws.Range("B4").Value = "."
$ws.Range("B5").Value = "."
$ws.Range("A6").Value = "."

# Update the active cell selection to match the recorded author session
$ws.Range("B14").Select()
